$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.281.97"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.666.82"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.42%  "

# Row 5: BNB
$ws.Range("D5").Value = "'218.47"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6: XRP
$ws.Range("E6").Value = "  -1.29%  "

# Row 7: USDC
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.39%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.2662"
$ws.Range("E8").Value = "  -1.19%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.06328"
$ws.Range("E9").Value = "  -1.04%  "

# Row 10: Solana
$ws.Range("D10").Value = "'20.97"
$ws.Range("E10").Value = "  -3.32%  "

# Row 11: TRON
$ws.Range("D11").Value = "'0.07759"
$ws.Range("E11").Value = "  -0.73%  "

# Row 12: WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.688.30"
$ws.Range("E12").Value = "  +0.99%  "

# Row 13: Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.444"
$ws.Range("E13").Value = "  -1.44%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.893.74"
$ws.Range("E14").Value = "  -0.67%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.5464"
$ws.Range("E15").Value = "  -1.71%  "

# Row 16: ShibaInu
$ws.Range("D16").Value = "0.0₅8249"
$ws.Range("E16").Value = "  -0.85%  "

# Row 17: Litecoin
$ws.Range("D17").Value = "'65.03"
$ws.Range("E17").Value = "  -0.86%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "26.300.29"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19: Dai
$ws.Range("E19").Value = "  +0.45%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'4.664"
$ws.Range("E20").Value = "  -1.53%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'194.89"
$ws.Range("E21").Value = "  +0.76%  "

# Row 22: Avalanche
$ws.Range("D22").Value = "'10.15"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23: Chainlink
$ws.Range("D23").Value = "'6.072"
$ws.Range("E23").Value = "  -4.27%  "

# Row 24: BinanceUSD
$ws.Range("D24").Value = "'1.008"
$ws.Range("E24").Value = "  +0.58%  "

# Row 25: Monero
$ws.Range("D25").Value = "'139.85"
$ws.Range("E25").Value = "  -1.79%  "

# Row 26: Stellar
$ws.Range("D26").Value = "'0.1240"
$ws.Range("E26").Value = "  -3.46%  "

# Row 27: Cosmos
$ws.Range("D27").Value = "'7.197"
$ws.Range("E27").Value = "  -2.73%  "

# Row 28: EthereumClassic
$ws.Range("E28").Value = "  -0.25%  "

# Row 29: Toncoin
$ws.Range("D29").Value = "'1.414"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30: Hedera
$ws.Range("D30").Value = "'0.06181"
$ws.Range("E30").Value = "  -0.99%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'1.281"
$ws.Range("E31").Value = "  +0.59%  "

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = "'3.598"
$ws.Range("E32").Value = "  -0.16%  "

# Row 33: Filecoin
$ws.Range("D33").Value = "'3.297"
$ws.Range("E33").Value = "  -4.29%  "

# Row 34: LidoDAOToken
$ws.Range("D34").Value = "'1.632"
$ws.Range("E34").Value = "  -2.71%  "

# Row 35: ARBITRUM
$ws.Range("D35").Value = "'0.9738"
$ws.Range("E35").Value = "  -3.31%  "

# Row 36: HuobiToken
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.427"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37: MXToken
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.782"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38: ImmutableX
$ws.Range("D38").Value = "'0.5766"
$ws.Range("E38").Value = "  -5.83%  "

# Row 39: VeChain
$ws.Range("D39").Value = "'0.01607"
$ws.Range("E39").Value = "  -1.32%  "

# Row 40: FraxShare
$ws.Range("D40").Value = "'6.028"
$ws.Range("E40").Value = "  -2.06%  "

# Row 41: TrustWalletToken
$ws.Range("D41").Value = "'0.8571"
$ws.Range("E41").Value = "  -0.76%  "

# Row 42: PaxDollar
$ws.Range("E42").Value = "  +0.47%  "

# Row 43: Maker
$ws.Range("D43").Value = "1.021.63"
$ws.Range("E43").Value = "  -5.68%  "

# Row 44: Quant
$ws.Range("D44").Value = "'100.30"
$ws.Range("E44").Value = "  +0.13%  "

# Row 45: RocketPoolETH
$ws.Range("D45").Value = "1.808.79"
$ws.Range("E45").Value = "  -0.75%  "

# Row 46: Aave
$ws.Range("D46").Value = "'57.79"
$ws.Range("E46").Value = "  +1.21%  "

# Row 47: Frax
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = "  +1.14%  "

# Row 48: EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.087"
$ws.Range("E48").Value = "  -0.94%  "

# Row 49: Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05187"
$ws.Range("E49").Value = "  -0.37%  "

# Row 50: RenderToken
$ws.Range("D50").Value = "'1.480"
$ws.Range("E50").Value = "  +0.56%  "

# Row 51: Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4219"
$ws.Range("E51").Value = "  -0.43%  "
